function Set-TextValue {
    param($range, $text)
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $origStyle
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws.Range("D2") "57.298.30"
Set-TextValue $ws.Range("E2") "  +2.95%  "

Set-TextValue $ws.Range("D3") "3.068.97"
Set-TextValue $ws.Range("E3") "  +5.43%  "

Set-TextValue $ws.Range("E4") "  +0.10%  "

Set-TextValue $ws.Range("D5") "513.03"
Set-TextValue $ws.Range("E5") "  +2.89%  "

Set-TextValue $ws.Range("D6") "142.25"
Set-TextValue $ws.Range("E6") "  +7.71%  "

Set-TextValue $ws.Range("E7") "  +0.10%  "

Set-TextValue $ws.Range("E8") "  +3.25%  "

Set-TextValue $ws.Range("D9") "7.29"
Set-TextValue $ws.Range("E9") "  +2.17%  "

Set-TextValue $ws.Range("D10") "0.108"
Set-TextValue $ws.Range("E10") "  +4.18%  "

Set-TextValue $ws.Range("E11") "  +6.46%  "

Set-TextValue $ws.Range("D12") "3.597.93"
Set-TextValue $ws.Range("E12") "  +5.53%  "

Set-TextValue $ws.Range("E13") "  +2.96%  "

Set-TextValue $ws.Range("D14") "25.89"
Set-TextValue $ws.Range("E14") "  +1.25%  "

Set-TextValue $ws.Range("E15") "  +3.69%  "

Set-TextValue $ws.Range("D16") "57.403.99"
Set-TextValue $ws.Range("E16") "  +3.43%  "

Set-TextValue $ws.Range("D17") "3.064.97"
Set-TextValue $ws.Range("E17") "  +5.32%  "

Set-TextValue $ws.Range("D18") "6.06"
Set-TextValue $ws.Range("E18") "  +2.08%  "

Set-TextValue $ws.Range("D19") "12.97"
Set-TextValue $ws.Range("E19") "  +3.11%  "

Set-TextValue $ws.Range("D20") "8.16"
Set-TextValue $ws.Range("E20") "  +6.76%  "

Set-TextValue $ws.Range("D21") "335.78"
Set-TextValue $ws.Range("E21") "  +7.53%  "

Set-TextValue $ws.Range("D22") "0.999"
Set-TextValue $ws.Range("E22") "  -0.11%  "

Set-TextValue $ws.Range("E23") "  +3.04%  "

Set-TextValue $ws.Range("D24") "65.30"
Set-TextValue $ws.Range("E24") "  +4.08%  "

Set-TextValue $ws.Range("D25") "0.171"
Set-TextValue $ws.Range("E25") "  +6.69%  "

Set-TextValue $ws.Range("E26") "  +0.02%  "

Set-TextValue $ws.Range("D27") "0.0₃0934"
Set-TextValue $ws.Range("E27") "  +13.18%  "

Set-TextValue $ws.Range("D28") "6.43"
Set-TextValue $ws.Range("E28") "  +2.85%  "

Set-TextValue $ws.Range("D29") "7.12"
Set-TextValue $ws.Range("E29") "  +5.76%  "

Set-TextValue $ws.Range("E30") "  +3.45%  "

Set-TextValue $ws.Range("D31") "20.73"
Set-TextValue $ws.Range("E31") "  +4.92%  "

Set-TextValue $ws.Range("E32") "  +4.19%  "

Set-TextValue $ws.Range("D33") "154.54"
Set-TextValue $ws.Range("E33") "  +2.20%  "

Set-TextValue $ws.Range("D34") "4.52"
Set-TextValue $ws.Range("E34") "  +4.18%  "

Set-TextValue $ws.Range("D35") "5.87"
Set-TextValue $ws.Range("E35") "  +5.60%  "

Set-TextValue $ws.Range("D36") "25.96"
Set-TextValue $ws.Range("E36") "  +10.18%  "

Set-TextValue $ws.Range("E37") "  +5.15%  "

Set-TextValue $ws.Range("E38") "  +6.12%  "

Set-TextValue $ws.Range("D39") "3.105.28"
Set-TextValue $ws.Range("E39") "  +5.65%  "

Set-TextValue $ws.Range("D40") "36.83"
Set-TextValue $ws.Range("E40") "  +1.26%  "

Set-TextValue $ws.Range("D41") "3.84"
Set-TextValue $ws.Range("E41") "  +4.50%  "

Set-TextValue $ws.Range("B42") "FirstDigitalUSD"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Range("D42") "1.00"
Set-TextValue $ws.Range("E42") "  +0.10%  "

Set-TextValue $ws.Range("B43") "Mantle"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D43") "0.667"
Set-TextValue $ws.Range("E43") "  +4.91%  "

Set-TextValue $ws.Range("D44") "2.265.32"
Set-TextValue $ws.Range("E44") "  +7.92%  "

Set-TextValue $ws.Range("B45") "VeChain"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D45") "0.0251"
Set-TextValue $ws.Range("E45") "  +8.07%  "

Set-TextValue $ws.Range("B46") "Stacks"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D46") "1.38"
Set-TextValue $ws.Range("E46") "  +4.74%  "

Set-TextValue $ws.Range("D47") "0.951"
Set-TextValue $ws.Range("E47") "  +4.04%  "

Set-TextValue $ws.Range("D48") "20.32"
Set-TextValue $ws.Range("E48") "  +9.98%  "

Set-TextValue $ws.Range("D49") "5.86"
Set-TextValue $ws.Range("E49") "  -1.88%  "

Set-TextValue $ws.Range("E50") "  +4.44%  "

Set-TextValue $ws.Range("D51") "0.686"
Set-TextValue $ws.Range("E51") "  +7.57%  "

